# Apply updated crypto price/volume data scraped by the GitHub Actions job.
# Mirrors a fresh run of the 'Updated cryptos list' workflow: refreshes the
# Price (column D) and Volume(1h) (column E) cells for each coin row (2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    # Row 2: Bitcoin
    @{ Cell = "D2"; Value = "26.978.93" }
    @{ Cell = "E2"; Value = "  -0.01%  " }
    # Row 3: Ethereum
    @{ Cell = "D3"; Value = "1.827.82" }
    @{ Cell = "E3"; Value = "  +0.45%  " }
    # Row 4: TetherUSD
    @{ Cell = "D4"; Value = "'1.006" }
    @{ Cell = "E4"; Value = "  -0.39%  " }
    # Row 5: BNB
    @{ Cell = "D5"; Value = "'311.08" }
    @{ Cell = "E5"; Value = "  +0.46%  " }
    # Row 6: USDC
    @{ Cell = "D6"; Value = "'1.004" }
    @{ Cell = "E6"; Value = "  -0.44%  " }
    # Row 7: XRP
    @{ Cell = "E7"; Value = "  -0.20%  " }
    # Row 8: Cardano
    @{ Cell = "D8"; Value = "'0.3711" }
    @{ Cell = "E8"; Value = "  +1.92%  " }
    # Row 9: Dogecoin
    @{ Cell = "D9"; Value = "'0.07337" }
    @{ Cell = "E9"; Value = "  +0.71%  " }
    # Row 10: Polygon
    @{ Cell = "D10"; Value = "'0.8772" }
    @{ Cell = "E10"; Value = "  +1.37%  " }
    # Row 11: TRON
    @{ Cell = "D11"; Value = "'0.07886" }
    @{ Cell = "E11"; Value = "  +3.49%  " }
    # Row 12: Solana
    @{ Cell = "D12"; Value = "'19.75" }
    @{ Cell = "E12"; Value = "  -0.39%  " }
    # Row 13: WrappedEther
    @{ Cell = "D13"; Value = "1.855.99" }
    @{ Cell = "E13"; Value = "  +1.16%  " }
    # Row 14: Polkadot
    @{ Cell = "D14"; Value = "'5.333" }
    @{ Cell = "E14"; Value = "  +0.07%  " }
    # Row 15: Chainlink
    @{ Cell = "D15"; Value = "'6.539" }
    @{ Cell = "E15"; Value = "  +1.08%  " }
    # Row 16: Litecoin
    @{ Cell = "D16"; Value = "'91.39" }
    @{ Cell = "E16"; Value = "  -1.85%  " }
    # Row 17: BinanceUSD
    @{ Cell = "D17"; Value = "'1.007" }
    @{ Cell = "E17"; Value = "  -0.19%  " }
    # Row 18: ShibaInu
    @{ Cell = "D18"; Value = "'0.000008841" }
    # Row 19: Dai
    @{ Cell = "E19"; Value = "  -0.54%  " }
    # Row 20: Avalanche
    @{ Cell = "D20"; Value = "'14.77" }
    @{ Cell = "E20"; Value = "  +2.07%  " }
    # Row 21: WrappedBTC
    @{ Cell = "D21"; Value = "27.018.24" }
    @{ Cell = "E21"; Value = "  -1.41%  " }
    # Row 22: Uniswap
    @{ Cell = "D22"; Value = "'5.102" }
    @{ Cell = "E22"; Value = "  -1.08%  " }
    # Row 23: Cosmos
    @{ Cell = "D23"; Value = "'10.52" }
    @{ Cell = "E23"; Value = "  -0.55%  " }
    # Row 24: WrappedliquidstakedEther2.0
    @{ Cell = "D24"; Value = "2.015.03" }
    @{ Cell = "E24"; Value = "  -4.44%  " }
    # Row 25: Monero
    @{ Cell = "D25"; Value = "'152.34" }
    @{ Cell = "E25"; Value = "  +0.35%  " }
    # Row 26: Toncoin
    @{ Cell = "D26"; Value = "'1.850" }
    @{ Cell = "E26"; Value = "  -0.33%  " }
    # Row 27: EthereumClassic
    @{ Cell = "E27"; Value = "  +1.72%  " }
    # Row 28: LidoDAOToken
    @{ Cell = "D28"; Value = "'2.036" }
    @{ Cell = "E28"; Value = "  -3.14%  " }
    # Row 29: InternetComputer(DFINITY)
    @{ Cell = "D29"; Value = "'5.109" }
    @{ Cell = "E29"; Value = "  +0.54%  " }
    # Row 30: BitcoinCash
    @{ Cell = "D30"; Value = "'115.64" }
    @{ Cell = "E30"; Value = "  -0.27%  " }
    # Row 31: Stellar
    @{ Cell = "D31"; Value = "'0.08891" }
    @{ Cell = "E31"; Value = "  -0.08%  " }
    # Row 32: HuobiToken
    @{ Cell = "D32"; Value = "'2.958" }
    @{ Cell = "E32"; Value = "  +0.15%  " }
    # Row 33: ImmutableX
    @{ Cell = "D33"; Value = "'0.7313" }
    @{ Cell = "E33"; Value = "  +0.04%  " }
    # Row 34: Filecoin
    @{ Cell = "D34"; Value = "'4.442" }
    @{ Cell = "E34"; Value = "  +0.28%  " }
    # Row 35: ARBITRUM
    @{ Cell = "D35"; Value = "'1.131" }
    @{ Cell = "E35"; Value = "  -0.70%  " }
    # Row 36: RenderToken
    @{ Cell = "D36"; Value = "'2.455" }
    @{ Cell = "E36"; Value = "  -2.76%  " }
    # Row 37: TrustWalletToken
    @{ Cell = "D37"; Value = "'1.077" }
    @{ Cell = "E37"; Value = "  +0.15%  " }
    # Row 38: VeChain
    @{ Cell = "D38"; Value = "'0.01948" }
    @{ Cell = "E38"; Value = "  +1.67%  " }
    # Row 39: Hedera
    @{ Cell = "D39"; Value = "'0.05217" }
    @{ Cell = "E39"; Value = "  -0.89%  " }
    # Row 40: MXToken
    @{ Cell = "D40"; Value = "'2.959" }
    @{ Cell = "E40"; Value = "  +0.83%  " }
    # Row 41: FraxShare
    @{ Cell = "D41"; Value = "'7.104" }
    @{ Cell = "E41"; Value = "  -0.49%  " }
    # Row 42: TheSandbox
    @{ Cell = "D42"; Value = "'0.5170" }
    @{ Cell = "E42"; Value = "  -0.88%  " }
    # Row 43: Algorand
    @{ Cell = "D43"; Value = "'0.1628" }
    @{ Cell = "E43"; Value = "  -0.22%  " }
    # Row 44: Aptos
    @{ Cell = "D44"; Value = "'8.155" }
    @{ Cell = "E44"; Value = "  -1.02%  " }
    # Row 45: Decentraland
    @{ Cell = "D45"; Value = "'0.4827" }
    @{ Cell = "E45"; Value = "  -0.59%  " }
    # Row 46: PaxDollar
    @{ Cell = "D46"; Value = "'1.004" }
    @{ Cell = "E46"; Value = "  -0.49%  " }
    # Row 47: EnergySwap
    @{ Cell = "D47"; Value = "'10.15" }
    @{ Cell = "E47"; Value = "  +0.44%  " }
    # Row 48: Quant
    @{ Cell = "D48"; Value = "'102.01" }
    @{ Cell = "E48"; Value = "  -1.15%  " }
    # Row 49: NEARProtocol
    @{ Cell = "D49"; Value = "'1.625" }
    @{ Cell = "E49"; Value = "  -0.51%  " }
    # Row 50: Cronos
    @{ Cell = "D50"; Value = "'0.06203" }
    @{ Cell = "E50"; Value = "  -0.33%  " }
    # Row 51: Aave
    @{ Cell = "D51"; Value = "'64.85" }
    @{ Cell = "E51"; Value = "  +0.47%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
